$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.636999999999999
$ws.Range("D2").Value = -7.33
$ws.Range("A3").Value = -21.535
$ws.Range("C3").Value = -12.505
$ws.Range("D6").Value = -7.958
$ws.Range("C12").Value = -11.536
$ws.Range("A14").Value = -21.682
$ws.Range("A16").Value = -21.373
$ws.Range("B18").Value = 5.414
$ws.Range("D19").Value = -8.078999999999999
$ws.Range("A21").Value = -20.524
$ws.Range("A23").Value = -20.703
$ws.Range("B24").Value = 6.915000000000001
$ws.Range("C24").Value = -13.492
$ws.Range("D24").Value = -7.509
$ws.Range("A25").Value = -20.386
$ws.Range("B25").Value = 6.748
$ws.Range("C25").Value = -12.701
$ws.Range("A26").Value = -21.458
$ws.Range("B27").Value = 6.443000000000001
$ws.Range("D27").Value = -8.327000000000002
$ws.Range("A29").Value = -21.15
$ws.Range("B30").Value = 5.852
$ws.Range("D30").Value = -7.648000000000001
$ws.Range("B31").Value = 6.306
$ws.Range("D31").Value = -8.193000000000001
$ws.Range("D33").Value = -7.904000000000001
$ws.Range("B39").Value = 7.811000000000002
$ws.Range("A40").Value = -20.294
$ws.Range("C41").Value = -12.61
$ws.Range("B42").Value = 8.395
$ws.Range("D42").Value = -8.472000000000001
$ws.Range("B48").Value = 5.156
$ws.Range("C50").Value = -13.087
$ws.Range("B51").Value = 5.216
$ws.Range("B52").Value = 4.933
$ws.Range("A53").Value = -21.793
$ws.Range("C53").Value = -11.857
$ws.Range("B55").Value = 4.581999999999999
$ws.Range("D55").Value = -8.086000000000002
$ws.Range("B56").Value = 6.069000000000001
$ws.Range("C56").Value = -12.88
$ws.Range("A57").Value = -21.322
$ws.Range("B57").Value = 6.104999999999999
$ws.Range("C57").Value = -12.872
$ws.Range("C58").Value = -12.826
$ws.Range("D58").Value = -8.094999999999999
$ws.Range("A59").Value = -21.852
$ws.Range("B60").Value = 5.005000000000001
$ws.Range("C61").Value = -13.181
$ws.Range("C63").Value = -11.757
$ws.Range("C64").Value = -11.805
$ws.Range("A65").Value = -21.521
$ws.Range("D65").Value = -7.885
$ws.Range("A69").Value = -21.543
$ws.Range("C70").Value = -11.621
$ws.Range("D70").Value = -7.571
$ws.Range("C72").Value = -11.799
$ws.Range("B73").Value = 6.186999999999999
$ws.Range("B74").Value = 8.825999999999999
$ws.Range("D74").Value = -8.049000000000001
$ws.Range("D75").Value = -7.711
$ws.Range("A79").Value = -21.18
$ws.Range("A83").Value = -21.372
$ws.Range("D83").Value = -8.526
$ws.Range("D84").Value = -7.934
$ws.Range("C86").Value = -12.495
$ws.Range("D86").Value = -7.705
$ws.Range("B89").Value = 5.938
$ws.Range("C89").Value = -11.222
$ws.Range("B90").Value = 5.5
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 6.027
$ws.Range("A93").Value = -21.324
$ws.Range("D96").Value = -7.418000000000001
$ws.Range("D97").Value = -8.122000000000002
$ws.Range("C98").Value = -12.712
$ws.Range("A100").Value = -21.61
$ws.Range("C100").Value = -12.887
$ws.Range("C102").Value = -12.801
